$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row of data (row 4): renewal / "Using Oldest Entry Date" case.
$ws.Range("A4").Value = "AAA_CSA"
$ws.Range("B4").Value = "SELECT"
$ws.Range("C4").Value = "CA"
$ws.Range("D4").Value = 20000102
$ws.Range("E4").Value = 20300102
$ws.Range("F4").Value = "SYMBOL_2000_CA_SELECT_ENTRY_DATE"

# A4:E4 pick up the same look as the rows above (style carries over as you
# extend the list); F4 keeps the formatting it already had before typing.
[void]$ws.Range("A3:E3").Copy()
[void]$ws.Range("A4:E4").PasteSpecial(-4122) # xlPasteFormats
$ws.Application.CutCopyMode = $false

# Widen column F to fit the new, longer symbol text.
[void]$ws.Columns.Item(6).Select()
[void]$ws.Columns.Item(6).EntireColumn.AutoFit()
